$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRbQ")

# Insert a new column before column B (shifts B..AF to C..AG),
# inheriting formatting from the column to the left (A).
$ws.Columns("B:B").Insert()

# New column B header: 2019 (style already inherited as s="4" from A1)
$ws.Range("B1").Value = 2019

# Give the new B data cells (rows 2-17) the same "integer, bold" style
# used throughout the rest of the header/body (style index 3).
$ws.Range("B2:B17").NumberFormat = "0"
$ws.Range("B2:B17").Font.Bold = $true

# Fill new column B with "= (value from the year to the right)" formulas.
$ws.Range("B2").Formula = "=C2"
$ws.Range("B3:B17").Formula = "=C3"

# Restore column B width to match column A (25.5703125 chars), since
# Excel's column insert otherwise leaves it at a default/best-fit width.
$ws.Columns("B:B").ColumnWidth = 24.67

# Update the active selection to match the authored workbook state.
$ws.Range("B2:B17").Select()
